$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 255, shifting rows 255:324 down to 256:325.
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new record's data.
$ws.Cells.Item(255, 1).Value = 5
$ws.Cells.Item(255, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(255, 3).Value = "Maule"
$ws.Cells.Item(255, 4).Value = 44932
$ws.Cells.Item(255, 5).Value = 7
$ws.Cells.Item(255, 6).Value = "Fruta"
$ws.Cells.Item(255, 7).Value = 100108
$ws.Cells.Item(255, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(255, 9).Value = 100108005
$ws.Cells.Item(255, 10).Value = "Piña"
$ws.Cells.Item(255, 11).Value = "Caramelo"
$ws.Cells.Item(255, 12).Value = "Segunda"
$ws.Cells.Item(255, 13).Value = 200
$ws.Cells.Item(255, 14).Value = 17000
$ws.Cells.Item(255, 15).Value = 17000
$ws.Cells.Item(255, 16).Value = 17000
$ws.Cells.Item(255, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(255, 18).Value = "Ecuador"
$ws.Cells.Item(255, 19).Value = 1214
$ws.Cells.Item(255, 20).Value = 14
